# ANN Lifetime Predictions - "updated models and predictions"
# Adds the 2023-12-19 slate of games to the Predictions sheet, refreshes the
# Pivot, and restores the window/selection state to the Predictions tab.

$wb = $excel.ActiveWorkbook
$wsPred = $wb.Worksheets.Item("Predictions")
$wsPivot = $wb.Worksheets.Item("Pivot")

# ---------------------------------------------------------------------------
# 1. Append the new game predictions (rows 12-22) to the Predictions sheet.
# ---------------------------------------------------------------------------
$newGames = @(
    @{ Row = 12; Winner = "Tampa Bay Lightning"; WinP = 0.81533999999999995; Loser = "St. Louis Blues";       LoseP = 0.18465999999999999; Site = "Playing At:  Tampa Bay Lightning   Home" },
    @{ Row = 13; Winner = "Boston Bruins";        WinP = 0.76312000000000002; Loser = "Minnesota Wild";        LoseP = 0.23688000000000001; Site = "Playing At:  Boston Bruins   Home" },
    @{ Row = 14; Winner = "Los Angeles Kings";    WinP = 0.76288999999999996; Loser = "San Jose Sharks";       LoseP = 0.23710999999999999; Site = "Playing At:  San Jose Sharks   Home" },
    @{ Row = 15; Winner = "Buffalo Sabres";       WinP = 0.72492999999999996; Loser = "Columbus Blue Jackets"; LoseP = 0.27506999999999998; Site = "Playing At:  Buffalo Sabres   Home" },
    @{ Row = 16; Winner = "New Jersey Devils";    WinP = 0.68744000000000005; Loser = "Philadelphia Flyers";   LoseP = 0.31256;              Site = "Playing At:  New Jersey Devils   Home" },
    @{ Row = 17; Winner = "Colorado Avalanche";   WinP = 0.67522000000000004; Loser = "Chicago Blackhawks";    LoseP = 0.32478000000000001; Site = "Playing At:  Chicago Blackhawks   Home" },
    @{ Row = 18; Winner = "Carolina Hurricanes";  WinP = 0.67174999999999996; Loser = "Vegas Golden Knights";  LoseP = 0.32824999999999999; Site = "Playing At:  Carolina Hurricanes   Home" },
    @{ Row = 19; Winner = "Edmonton Oilers";      WinP = 0.64061000000000001; Loser = "New York Islanders";    LoseP = 0.35938999999999999; Site = "Playing At:  New York Islanders   Home" },
    @{ Row = 20; Winner = "Ottawa Senators";      WinP = 0.58294999999999997; Loser = "Arizona Coyotes";       LoseP = 0.41704999999999998; Site = "Playing At:  Arizona Coyotes   Home" },
    @{ Row = 21; Winner = "Nashville Predators";  WinP = 0.58067000000000002; Loser = "Vancouver Canucks";     LoseP = 0.41932999999999998; Site = "Playing At:  Nashville Predators   Home" },
    @{ Row = 22; Winner = "Toronto Maple Leafs";  WinP = 0.52131000000000005; Loser = "New York Rangers";      LoseP = 0.47869;              Site = "Playing At:  Toronto Maple Leafs   Home" }
)

foreach ($g in $newGames) {
    $r = $g.Row

    $wsPred.Range("A$r").Value = 45279
    $wsPred.Range("B$r").Value = $g.Winner
    $wsPred.Range("C$r").Value = $g.WinP
    $wsPred.Range("D$r").Value = $g.Loser
    $wsPred.Range("E$r").Value = $g.LoseP
    $wsPred.Range("F$r").Value = $g.Site

    $formula = '=_xlfn.IFS(C' + $r + ' >= 0.85, "85  <", AND(C' + $r + ' >=0.8, C' + $r + ' < 0.85), "80-85", AND(C' + $r + ' >= 0.7, C' + $r + ' < 0.8), "70-80", AND(C' + $r + ' >= 0.6, C' + $r + ' < 0.7),  "60-70", AND(C' + $r + ' >= 0.5, C' + $r + ' < 0.6), "50-60")'
    $wsPred.Range("H$r").FormulaArray = $formula
}

# ---------------------------------------------------------------------------
# 2. Refresh the pivot table / pivot cache so it reflects the new source rows.
# ---------------------------------------------------------------------------
$pt = $wsPivot.PivotTables(1)
$pt.PivotCache().Refresh()

# ---------------------------------------------------------------------------
# 3. Restore the Pivot sheet's B4:C7 number-format cells back to General
#    (matches the author clearing the stray "applyNumberFormat" style).
# ---------------------------------------------------------------------------
$wsPivot.Range("B4:C7").ClearFormats()

# ---------------------------------------------------------------------------
# 4. Window / selection state: Predictions becomes the active (visible) tab,
#    with E16 selected; Pivot keeps its prior B12 selection.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Left = 1740
$win.Top = 1185
$win.Width = 19125
$win.Height = 11235

$wsPivot.Range("B12").Select()
$wsPred.Activate()
$wsPred.Range("E16").Select()
